$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37/38: coin swap - VeChain/Hedera order flipped with new values
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '''0.06111'
$ws.Range("E37").Value = '  +2.88%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.02264'
$ws.Range("E38").Value = '  +2.48%  '

# Price (D) and Volume/1h (E) updates for remaining rows
$ws.Range("D2").Value = '27.904.36'
$ws.Range("E2").Value = '  +2.14%  '
$ws.Range("D3").Value = '1.880.96'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''333.11'
$ws.Range("E5").Value = '  +3.73%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  +5.88%  '
$ws.Range("E8").Value = '  +4.27%  '
$ws.Range("D9").Value = '''48.47'
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").Value = '''0.08056'
$ws.Range("E10").Value = '  +2.79%  '
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").Value = '''21.96'
$ws.Range("E12").Value = '  +3.07%  '
$ws.Range("D13").Value = '1.906.13'
$ws.Range("E13").Value = '  +3.20%  '
$ws.Range("D14").Value = '''5.970'
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").Value = '''7.193'
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").Value = '''1.006'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").Value = '''0.00001052'
$ws.Range("E17").Value = '  +2.23%  '
$ws.Range("D18").Value = '''87.34'
$ws.Range("E18").Value = '  +2.07%  '
$ws.Range("D19").Value = '''0.06622'
$ws.Range("E19").Value = '  +1.91%  '
$ws.Range("D20").Value = '''17.40'
$ws.Range("E20").Value = '  +3.02%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '28.076.35'
$ws.Range("E22").Value = '  +2.83%  '
$ws.Range("D23").Value = '''5.511'
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("D24").Value = '''11.08'
$ws.Range("E24").Value = '  +2.83%  '
$ws.Range("D25").Value = '''2.317'
$ws.Range("E25").Value = '  +3.05%  '
$ws.Range("D26").Value = '2.137.55'
$ws.Range("E26").Value = '  +3.57%  '
$ws.Range("D27").Value = '''157.22'
$ws.Range("E27").Value = '  +3.90%  '
$ws.Range("D28").Value = '''20.27'
$ws.Range("E28").Value = '  +4.92%  '
$ws.Range("D29").Value = '''2.109'
$ws.Range("E29").Value = '  +2.95%  '
$ws.Range("D30").Value = '''5.633'
$ws.Range("E30").Value = '  +2.37%  '
$ws.Range("D31").Value = '''122.75'
$ws.Range("E31").Value = '  +2.98%  '
$ws.Range("D32").Value = '''0.9824'
$ws.Range("E32").Value = '  +5.99%  '
$ws.Range("D33").Value = '''0.09569'
$ws.Range("E33").Value = '  +2.82%  '
$ws.Range("D34").Value = '''1.466'
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").Value = '''3.633'
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").Value = '''5.331'
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("D40").Value = '''8.253'
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").Value = '''0.6039'
$ws.Range("E41").Value = '  +2.74%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").Value = '''0.1912'
$ws.Range("E43").Value = '  +3.59%  '
$ws.Range("D44").Value = '''10.34'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("D45").Value = '''0.5746'
$ws.Range("E45").Value = '  +2.23%  '
$ws.Range("D46").Value = '''1.249'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '''12.34'
$ws.Range("E47").Value = '  +2.07%  '
$ws.Range("D48").Value = '''3.416'
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("D49").Value = '''1.951'
$ws.Range("E49").Value = '  +1.87%  '
$ws.Range("D50").Value = '''0.06833'
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").Value = '''113.97'
$ws.Range("E51").Value = '  +5.35%  '
